$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-22 Tuesday", "2025-07-23 Wednesday"),
    @("115×7=805", "854×7=5978"),
    @("968×5=4840", "779×8=6232"),
    @("332×5=1660", "188×4=752"),
    @("555×9=4995", "283×6=1698"),
    @("916×6=5496", "960×8=7680"),
    @("363×5=1815", "823×4=3292"),
    @("809×8=6472", "489×2=978"),
    @("956×8=7648", "604×3=1812"),
    @("577×5=2885", "721×4=2884"),
    @("921×7=6447", "628×3=1884"),
    @("108×6=648", "581×5=2905"),
    @("276×4=1104", "539×4=2156"),
    @("242×7=1694", "643×4=2572"),
    @("605×7=4235", "867×5=4335"),
    @("328×3=984", "746×8=5968"),
    @("655×3=1965", "690×8=5520"),
    @("660×3=1980", "552×4=2208"),
    @("142×8=1136", "495×9=4455"),
    @("866×9=7794", "482×5=2410"),
    @("603×4=2412", "807×9=7263"),
    @("628×5=3140", "523×3=1569"),
    @("925×5=4625", "131×8=1048"),
    @("973×2=1946", "182×7=1274"),
    @("232×2=464", "630×8=5040"),
    @("439×6=2634", "492×5=2460")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
